# Weekly price update: insert a new week's record as row 26
# (Vega Monumental Concepción - Pepino dulce, "Primera" quality),
# pushing the existing rows 26-64 down to 27-65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 26..64 down to 27..65, leaving a blank row 26 that
# inherits the formatting (incl. the date style) of the row above it.
$ws.Rows("26:26").Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A26").Value = 11
$ws.Range("B26").Value = "Vega Monumental Concepción"
$ws.Range("C26").Value = "Bíobío"
$ws.Range("D26").Value = 45028
$ws.Range("E26").Value = 8
$ws.Range("F26").Value = 100112043
$ws.Range("G26").Value = "Pepino dulce"
$ws.Range("H26").Value = "Cultivar IV Región"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 220
$ws.Range("K26").Value = 13000
$ws.Range("L26").Value = 14000
$ws.Range("M26").Value = 13455
$ws.Range("N26").Value = "$/bandeja 18 kilos"
$ws.Range("O26").Value = "Provincia de Limarí"
$ws.Range("P26").Value = 748
$ws.Range("Q26").Value = 18
$ws.Range("R26").Value = "Hortaliza"
